$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Validation")

$ws.Range("A3").Value = 178
$ws.Range("B3").Value = 210
$ws.Range("C3").Value = 246
$ws.Range("D3").Value = 1282
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 11
